$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Step 1: Clear all of the "body" content (paragraphs 3..13), merging it down
# to a single, cleanly-formatted empty paragraph so we can rebuild it without
# carrying over any stray bold/italic paragraph-mark formatting.
# ---------------------------------------------------------------------------
$startRng = $d.Paragraphs.Item(3).Range.Start
$endRng = $d.Paragraphs.Item(13).Range.End
$clearRng = $d.Range($startRng, $endRng)
$clearRng.Delete()

# ---------------------------------------------------------------------------
# Step 2: Rebuild the body paragraphs, left-aligned, with plain (non-bold,
# non-italic) runs except for the three section headings.
# ---------------------------------------------------------------------------

# Paragraph 3: blank spacer line
$p3 = $d.Paragraphs.Item(3)
$p3.Format.Alignment = 0

$p3.Range.InsertParagraphAfter()
# Paragraph 4: "What went right?" heading
$p4 = $d.Paragraphs.Item(4)
$p4.Format.Alignment = 0
$p4.Range.InsertAfter("What went right?")

$p4.Range.InsertParagraphAfter()
# Paragraph 5: body text describing what went right
$p5 = $d.Paragraphs.Item(5)
$p5.Format.Alignment = 0
$p5.Range.InsertAfter("As two of the platforms use very similar controls and hardware, I only needed to create one build for both Oculus Rift S and Oculus Quest 2. Unity automatically down-sampled the textures and lighting effects for the Oculus build once the build target changed, saving development time by pre-emptively optimizing for the weaker hardware. All players reported that the game ran smoothly on their respective platforms and hardware. ")

$p5.Range.InsertParagraphAfter()
# Paragraph 6: blank spacer line
$p6 = $d.Paragraphs.Item(6)
$p6.Format.Alignment = 0

$p6.Range.InsertParagraphAfter()
# Paragraph 7: "What went wrong?" heading
$p7 = $d.Paragraphs.Item(7)
$p7.Format.Alignment = 0
$p7.Range.InsertAfter("What went wrong?")

$p7.Range.InsertParagraphAfter()
# Paragraph 8: body text describing what went wrong
$p8 = $d.Paragraphs.Item(8)
$p8.Format.Alignment = 0
$p8.Range.InsertAfter("Building an Android version of the game for Oculus headsets was much more complex than building for the HP Reverb and other Windows platforms. I had to download additional components to build the game. When building the game, I discovered bugs that could not happen in the development setting, which occurred when playing the game more than once. These issues took additional development time to resolve.")

$p8.Range.InsertParagraphAfter()
# Paragraph 9: blank spacer line
$p9 = $d.Paragraphs.Item(9)
$p9.Format.Alignment = 0

$p9.Range.InsertParagraphAfter()
# Paragraph 10: "What can be improved for next time?" heading
$p10 = $d.Paragraphs.Item(10)
$p10.Format.Alignment = 0
$p10.Range.InsertAfter("What can be improved for next time?")

$p10.Range.InsertParagraphAfter()
# Paragraph 11: body text describing what can be improved
$p11 = $d.Paragraphs.Item(11)
$p11.Format.Alignment = 0
$p11.Range.InsertAfter("For future multiplatform releases, it would be very beneficial to test the game on the various target platforms sooner in the development cycle to assist in identifying these platform-specific issues so they can be addressed quicker and in a more thought-out way. One other improvement would be to check which platforms were available to those who would be testing my game, and target those specifically to streamline the testing process.")

# ---------------------------------------------------------------------------
# Step 3: Bold the three section headings (text only, so the new body runs
# stay completely unformatted).
# ---------------------------------------------------------------------------
$headings = @("What went right?", "What went wrong?", "What can be improved for next time?")
foreach ($h in $headings) {
    $findRng = $d.Content
    $findRng.Find.Execute($h, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if ($findRng.Find.Found) {
        $findRng.Font.Bold = 1
        $findRng.Font.BoldBi = 1
    }
}

# ---------------------------------------------------------------------------
# Step 4: Mark the "down-sampled" word with a bookmark, matching the source
# edit (left over from an accepted inline rewrite suggestion).
# ---------------------------------------------------------------------------
$bmRng = $d.Content
$bmRng.Find.Execute("down-sampled", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($bmRng.Find.Found) {
    $d.Bookmarks.Add("_Int_Kmklhv0B", $bmRng)
}

Write-Output ("Final paragraph count: " + $d.Paragraphs.Count)
